# Add a new reservation row (row 13) to the "Order Data" sheet,
# mirroring the existing reservation rows (e.g. row 12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 13

$ws.Cells.Item($row, 1).Value = 12          # Reservation ID
$ws.Cells.Item($row, 2).Value = 2           # Number of Guests
$ws.Cells.Item($row, 3).Value = 45412       # Date Of Reservation (serial date)
$ws.Cells.Item($row, 3).NumberFormat = "yyyy-MM-dd"
$ws.Cells.Item($row, 4).Value = 4           # Customer ID
$ws.Cells.Item($row, 5).Value = "2 seat"    # Table Type
$ws.Cells.Item($row, 6).Value = 1           # Table Count
